$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with rich-text runs; flattened to plain text, same font) ---
$ws.Range("A8").Value = "Volume 32   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/10/2025  Through  2/16/2025"

# --- Cells changing from text placeholder ("-" / "***.*") to numeric: set NumberFormat first ---
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("C16").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 2
$ws.Range("C31").NumberFormat = "#,##0"
$ws.Range("C31").Value = 2
$ws.Range("F31").NumberFormat = "#,##0"
$ws.Range("F31").Value = 2
$ws.Range("I31").NumberFormat = "#,##0"
$ws.Range("I31").Value = 2
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E22").Value = -50

# --- Plain numeric value updates (style/number-format unchanged) ---
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 14
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = 40
$ws.Range("L16").Value = 55.555555555555
$ws.Range("M16").Value = -41.666666666666
$ws.Range("N16").Value = -87.272727272727
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 15
$ws.Range("H17").Value = 150
$ws.Range("I17").Value = 21
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = 90.909090909090
$ws.Range("L17").Value = 133.333333333333
$ws.Range("M17").Value = 133.333333333333
$ws.Range("N17").Value = 75
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 20
$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 12
$ws.Range("K18").Value = 166.666666666667
$ws.Range("L18").Value = 23.076923076923
$ws.Range("M18").Value = 128.571428571429
$ws.Range("N18").Value = -84.236453201970
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -17.857142857142
$ws.Range("J19").Value = 45
$ws.Range("K19").Value = -4.444444444444
$ws.Range("L19").Value = -17.307692307692
$ws.Range("M19").Value = -25.862068965517
$ws.Range("N19").Value = -61.946902654867
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -30
$ws.Range("I20").Value = 18
$ws.Range("J20").Value = 12
$ws.Range("K20").Value = 50
$ws.Range("L20").Value = -10
$ws.Range("M20").Value = 12.5
$ws.Range("N20").Value = -96.571428571428
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = -15.384615384615
$ws.Range("F21").Value = 73
$ws.Range("G21").Value = 56
$ws.Range("H21").Value = 30.357142857142
$ws.Range("I21").Value = 130
$ws.Range("J21").Value = 91
$ws.Range("K21").Value = 42.857142857142
$ws.Range("L21").Value = 11.111111111111
$ws.Range("M21").Value = 7.438016528925
$ws.Range("N21").Value = -86.528497409326
$ws.Range("C22").Value = 1
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = 25
$ws.Range("L22").Value = 66.666666666666
$ws.Range("M22").Value = 25
$ws.Range("C24").Value = 55
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = 44.736842105263
$ws.Range("F24").Value = 178
$ws.Range("G24").Value = 134
$ws.Range("H24").Value = 32.835820895522
$ws.Range("I24").Value = 274
$ws.Range("J24").Value = 208
$ws.Range("K24").Value = 31.730769230769
$ws.Range("L24").Value = 26.851851851851
$ws.Range("M24").Value = 194.623655913978
$ws.Range("C25").Value = 47
$ws.Range("D25").Value = 35
$ws.Range("E25").Value = 34.285714285714
$ws.Range("F25").Value = 149
$ws.Range("G25").Value = 105
$ws.Range("H25").Value = 41.904761904761
$ws.Range("I25").Value = 227
$ws.Range("J25").Value = 163
$ws.Range("K25").Value = 39.263803680981
$ws.Range("L25").Value = 43.670886075949
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 150
$ws.Range("F26").Value = 33
$ws.Range("H26").Value = 135.714285714286
$ws.Range("I26").Value = 47
$ws.Range("J26").Value = 23
$ws.Range("K26").Value = 104.347826086957
$ws.Range("L26").Value = 80.769230769230
$ws.Range("M26").Value = 51.612903225806
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 400
$ws.Range("I28").Value = 6
$ws.Range("K28").Value = 100
$ws.Range("L28").Value = 20
$ws.Range("L31").Value = 100
